# Generate Report for Handoff
# Updates status / timestamp cells to reflect a new "Ready for handoff" report
# generation pass, and widens the affected datetime columns to fit the new
# timestamp text.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Overview sheet ---------------------------------------------------
# Status columns for zh-cn (E2) and de-de (F2)
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
# Latest HO Xliff Generate Date
$wsOverview.Range("G2").Value = "2016-08-19 06:56:51"

# Widen the two status-adjacent datetime columns (E & F). The target
# character width from the source edit is 17.2159881591797, but the
# engine quantizes ColumnWidth to integer pixels (1/6-character steps
# at this font's 7px max-digit-width), so an input of 98/6 is the
# closest achievable round-trip (-> 17.166666666666668).
$targetColWidth = 98 / 6
$wsOverview.Columns.Item(5).ColumnWidth = $targetColWidth
$wsOverview.Columns.Item(6).ColumnWidth = $targetColWidth

# --- zh-cn sheet --------------------------------------------------------
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-08-19 06:56:46"
$wsZhCn.Columns.Item(3).ColumnWidth = $targetColWidth

# --- de-de sheet --------------------------------------------------------
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-08-19 06:56:51"
$wsDeDe.Columns.Item(3).ColumnWidth = $targetColWidth

$wb.Save()
